$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (pushes existing row 4+ down by one), scoped
# to the used columns (A:J) so unused row/column metadata isn't disturbed.
$ws.Range("A4:J4").Insert(-4121)  # xlShiftDown

# Copy the formatting of row 3 (the row above) into the newly inserted row 4,
# matching the "last row of block" bottom-border style it had before insertion.
$ws.Range("A3:H3").Copy()
$ws.Range("A4:H4").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row's values
$ws.Range("A4").Value = 0
$ws.Range("B4").Value = "Facility"
$ws.Range("C4").Value = "Global"
$ws.Range("D4").Value = "Vacuum chamber"
$ws.Range("E4").Value = "Mother volume radius"
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = "m"

# Row 3 above loses its bottom border (no longer last row of its block) while
# keeping the rest of its formatting; F3 additionally becomes right-aligned.
$ws.Range("A3:E3,G3:H3").Borders.Item(9).LineStyle = -4142  # xlEdgeBottom, xlLineStyleNone
$ws.Range("F3").Borders.Item(9).LineStyle = -4142
$ws.Range("F3").HorizontalAlignment = -4152  # xlRight

# Move selection as recorded in the edited workbook
$ws.Range("F4").Select()
